$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2775.3572
$ws.Range("I28").Value = 983.4
$ws.Range("K28").Value = 983.4
$ws.Range("M28").Value = -498.4
$ws.Range("H61").Value = 2033.8
$ws.Range("I61").Value = 417.25
$ws.Range("J61").Value = 8500
$ws.Range("K61").Value = 1251.75
$ws.Range("L61").Value = 25500
$ws.Range("M61").Value = -1079.75
$ws.Range("N61").Value = -25844
$ws.Range("H86").Value = 2533.2307
$ws.Range("I86").Value = 2193.6667
$ws.Range("J86").Value = 2824.2856
$ws.Range("K86").Value = 2193.6667
$ws.Range("L86").Value = 2824.2856
$ws.Range("M86").Value = -1070.6667
$ws.Range("N86").Value = -5070.2856
$ws.Range("H89").Value = 2533.2307
$ws.Range("I89").Value = 2193.6667
$ws.Range("J89").Value = 2824.2856
$ws.Range("K89").Value = 10968.3335
$ws.Range("L89").Value = 14121.428
$ws.Range("M89").Value = -5352.333500000001
$ws.Range("N89").Value = -25353.428
$ws.Range("H98").Value = 1924.25
$ws.Range("I98").Value = 1088.8422
$ws.Range("K98").Value = 1088.8422
$ws.Range("M98").Value = 409.1578
$ws.Range("H112").Value = 2461.0667
$ws.Range("I112").Value = 1672.5
$ws.Range("J112").Value = 2582.3845
$ws.Range("K112").Value = 5017.5
$ws.Range("L112").Value = 7747.1535
$ws.Range("M112").Value = -3909.5
$ws.Range("N112").Value = -9963.1535
$ws.Range("H122").Value = 1924.25
$ws.Range("I122").Value = 1088.8422
$ws.Range("K122").Value = 3266.5266
$ws.Range("M122").Value = -816.5266000000001
$ws.Range("H129").Value = 1257.7222
$ws.Range("I129").Value = 1064.8334
$ws.Range("J129").Value = 1354.1666
$ws.Range("K129").Value = 3194.5002
$ws.Range("L129").Value = 4062.4998
$ws.Range("M129").Value = 1805.4998
$ws.Range("N129").Value = -14062.4998
$ws.Range("H132").Value = 6335.5713
$ws.Range("I132").Value = 4927.9355
$ws.Range("J132").Value = 17244.75
$ws.Range("K132").Value = 14783.8065
$ws.Range("L132").Value = 51734.25
$ws.Range("M132").Value = -12253.8065
$ws.Range("N132").Value = -56794.25
$ws.Range("H138").Value = 3072.1738
$ws.Range("J138").Value = 3115.5667
$ws.Range("L138").Value = 9346.7001
$ws.Range("N138").Value = -19626.7001
$ws.Range("H141").Value = 3098.3333
$ws.Range("I141").Value = 1889.174
$ws.Range("J141").Value = 7071.2856
$ws.Range("K141").Value = 5667.522
$ws.Range("L141").Value = 21213.8568
$ws.Range("M141").Value = -487.5219999999999
$ws.Range("N141").Value = -31573.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2573.16
$ws.Range("I2").Value = 1597.1052
$ws.Range("J2").Value = 5664
$ws.Range("K2").Value = 1597.1052
$ws.Range("L2").Value = 5664
$ws.Range("M2").Value = -1484.1052
$ws.Range("N2").Value = -5890
$ws.Range("H32").Value = 30024.736
$ws.Range("I32").Value = 43665.293
$ws.Range("K32").Value = 43665.293
$ws.Range("M32").Value = -43378.293
$ws.Range("H45").Value = 2044.0416
$ws.Range("I45").Value = 1336.8889
$ws.Range("K45").Value = 1336.8889
$ws.Range("M45").Value = -959.8888999999999
$ws.Range("H116").Value = 2573.16
$ws.Range("I116").Value = 1597.1052
$ws.Range("J116").Value = 5664
$ws.Range("K116").Value = 1597.1052
$ws.Range("L116").Value = 5664
$ws.Range("M116").Value = 696.8948
$ws.Range("N116").Value = -10252
$ws.Range("H117").Value = 149999
$ws.Range("J117").Value = 149999
$ws.Range("L117").Value = 149999
$ws.Range("N117").Value = -159177

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2573.16
$ws.Range("I3").Value = 1597.1052
$ws.Range("J3").Value = 5664
$ws.Range("K3").Value = 1597.1052
$ws.Range("L3").Value = 5664
$ws.Range("M3").Value = -1483.1052
$ws.Range("N3").Value = -5892
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H86").Value = 4099.3076
$ws.Range("I86").Value = 2948
$ws.Range("J86").Value = 4611
$ws.Range("K86").Value = 2948
$ws.Range("L86").Value = 4611
$ws.Range("M86").Value = -1825
$ws.Range("N86").Value = -6857
$ws.Range("H89").Value = 4099.3076
$ws.Range("I89").Value = 2948
$ws.Range("J89").Value = 4611
$ws.Range("K89").Value = 14740
$ws.Range("L89").Value = 23055
$ws.Range("M89").Value = -9124
$ws.Range("N89").Value = -34287
$ws.Range("H99").Value = 206792
$ws.Range("I99").Value = 502005
$ws.Range("K99").Value = 502005
$ws.Range("M99").Value = -500507
$ws.Range("H107").Value = 1139.6923
$ws.Range("I107").Value = 918
$ws.Range("K107").Value = 918
$ws.Range("M107").Value = 1002
$ws.Range("H112").Value = 120459.5
$ws.Range("J112").Value = 120459.5
$ws.Range("L112").Value = 120459.5
$ws.Range("N112").Value = -123413.5
$ws.Range("H137").Value = 69999
$ws.Range("J137").Value = 69999
$ws.Range("L137").Value = 69999
$ws.Range("N137").Value = -80199

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 10600
$ws.Range("J28").Value = 10600
$ws.Range("L28").Value = 10600
$ws.Range("N28").Value = -11090
$ws.Range("H86").Value = 4835.2856
$ws.Range("I86").Value = 4569.4
$ws.Range("K86").Value = 4569.4
$ws.Range("M86").Value = -3446.4
$ws.Range("H89").Value = 4835.2856
$ws.Range("I89").Value = 4569.4
$ws.Range("K89").Value = 22847
$ws.Range("M89").Value = -17231
$ws.Range("H99").Value = 2020.5555
$ws.Range("I99").Value = 1950
$ws.Range("J99").Value = 2077
$ws.Range("K99").Value = 1950
$ws.Range("L99").Value = 2077
$ws.Range("M99").Value = -452
$ws.Range("N99").Value = -5073
$ws.Range("H104").Value = 52094.668
$ws.Range("J104").Value = 68142
$ws.Range("L104").Value = 68142
$ws.Range("N104").Value = -73384
$ws.Range("H115").Value = 64998.5
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 64998.5
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 64998.5
$ws.Range("M115").ClearContents()
$ws.Range("N115").Value = -67348.5
$ws.Range("H126").Value = 2020.5555
$ws.Range("I126").Value = 1950
$ws.Range("J126").Value = 2077
$ws.Range("K126").Value = 5850
$ws.Range("L126").Value = 6231
$ws.Range("M126").Value = -3380
$ws.Range("N126").Value = -11171
$ws.Range("H132").Value = 2472.1667
$ws.Range("I132").Value = 2472.1667
$ws.Range("K132").Value = 7416.500100000001
$ws.Range("M132").Value = -4886.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 700529.5600000001
$ws.Range("I4").Value = 700529.5600000001
$ws.Range("K4").Value = 2101588.68
$ws.Range("M4").Value = -2101476.68
$ws.Range("H17").Value = 45.733334
$ws.Range("J17").Value = 42
$ws.Range("L17").Value = 126
$ws.Range("N17").Value = -464
$ws.Range("H93").Value = 9874.75
$ws.Range("J93").Value = 9874.75
$ws.Range("L93").Value = 29624.25
$ws.Range("N93").Value = -33368.25
$ws.Range("H122").Value = 1353.909
$ws.Range("I122").Value = 998.3333
$ws.Range("J122").Value = 1487.25
$ws.Range("K122").Value = 8984.9997
$ws.Range("L122").Value = 13385.25
$ws.Range("M122").Value = -6534.9997
$ws.Range("N122").Value = -18285.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 37000
$ws.Range("J101").Value = 37000
$ws.Range("L101").Value = 37000
$ws.Range("N101").Value = -43490
$ws.Range("H108").Value = 48999.8
$ws.Range("J108").Value = 48999.8
$ws.Range("L108").Value = 48999.8
$ws.Range("N108").Value = -56679.8
$ws.Range("H129").Value = 69999
$ws.Range("J129").Value = 69999
$ws.Range("L129").Value = 69999
$ws.Range("N129").Value = -79999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 2831.25
$ws.Range("I32").Value = 2831.25
$ws.Range("K32").Value = 2831.25
$ws.Range("M32").Value = -2514.25
$ws.Range("H40").Value = 5590.1
$ws.Range("I40").Value = 5266.778
$ws.Range("J40").Value = 8500
$ws.Range("K40").Value = 5266.778
$ws.Range("L40").Value = 8500
$ws.Range("M40").Value = -5130.778
$ws.Range("N40").Value = -8772
$ws.Range("H127").Value = 73333
$ws.Range("J127").Value = 73333
$ws.Range("L127").Value = 73333
$ws.Range("N127").Value = -83253
$ws.Range("H132").Value = 70549.5
$ws.Range("I132").Value = 241513.8
$ws.Range("J132").Value = 4794
$ws.Range("K132").Value = 724541.3999999999
$ws.Range("L132").Value = 14382
$ws.Range("M132").Value = -722011.3999999999
$ws.Range("N132").Value = -19442
$ws.Range("H133").Value = 85229.25
$ws.Range("J133").Value = 85229.25
$ws.Range("L133").Value = 85229.25
$ws.Range("N133").Value = -90289.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 114787.22
$ws.Range("I132").Value = 128510.625
$ws.Range("K132").Value = 385531.875
$ws.Range("M132").Value = -383001.875
$ws.Range("H133").Value = 79994.5
$ws.Range("J133").Value = 79994.5
$ws.Range("L133").Value = 79994.5
$ws.Range("N133").Value = -90114.5
